$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the old "Meta description: ..." paragraph that currently sits
#    right after the title heading ("Play Dancing in Rio Slot for Free -
#    Review"). It consists of an empty run, a bold "Meta description" run
#    and a plain run with the description text - deleting the whole
#    paragraph range (including its end-of-paragraph mark) removes it
#    cleanly.
# ---------------------------------------------------------------------------
$metaLabel = "Meta description"
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.StartsWith($metaLabel)) {
        $para.Range.Delete()
        $found = $true
        break
    }
}

# ---------------------------------------------------------------------------
# 2) Insert a new bold paragraph "Play Dancing in Rio Slot for Free -
#    Review" right before the final paragraph of the document (the one
#    that used to hold the image-generation prompt). Use InsertXML with a
#    minimal WordprocessingML package so the resulting markup matches the
#    structure used elsewhere in the document (leading empty run followed
#    by a bold run).
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertPoint = $d.Range($lastPara.Range.Start - 1, $lastPara.Range.Start - 1)

$xmlFragment = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
            '<w:p>' +
              '<w:r/>' +
              '<w:r><w:rPr><w:b/></w:rPr><w:t>Play Dancing in Rio Slot for Free - Review</w:t></w:r>' +
            '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

$null = $insertPoint.InsertXML($xmlFragment)

# ---------------------------------------------------------------------------
# 3) Replace the text of the (now) last paragraph - formerly the image
#    generation prompt, still italic - with the meta-description copy.
# ---------------------------------------------------------------------------
$oldPrompt = "Create a feature image for the Dancing in Rio game in a cartoon style. The image should feature a happy Maya warrior with glasses. The Maya warrior should be surrounded by other happy carnival dancers in brightly colored costumes, all dancing in the streets with confetti raining down on them. The image should capture the joyful and vibrant vibes of the Rio Carnival. The background should be the iconic Christ the Redeemer statue, adding a touch of authenticity to the image. The overall feel of the image should be fun and lively, inviting players to join in on the carnival festivities and try their luck at this exciting slot game."
$newMeta = "Read our review of Dancing in Rio to discover its features and see if you want to play it for free. Includes a progressive jackpot and up to 50 free spins."

$replaced = $d.Content.Find.Execute($oldPrompt, $true, $false, $false, $false, $false, `
                                     $true, 1, $false, $newMeta, 2)

if (-not $found) {
    Write-Output "WARNING: Meta description paragraph was not found for deletion."
}
if (-not $replaced) {
    Write-Output "WARNING: Image-prompt paragraph text was not found/replaced."
}
